$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest six years of data (2004年-2009年, rows 2-7).
# Remaining rows shift up so 2010年 becomes row 2, ..., 2020年 becomes row 12.
$ws.Range("A2:H7").Delete()

# Copy the formatting of the last existing data row (now row 12, 2020年) down
# to the new row 13 so the year label picks up the same style (s="1").
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Append the new 2021年 data row.
$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 1311.8146
$ws.Range("E13").Value = 488.989
$ws.Range("F13").Value = 375.3709
$ws.Range("H13").Value = 405.8842
